$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.502111673355103
$ws.Range("B1").Value = 5.106636524200439
$ws.Range("C1").Value = 0.4210483431816101
$ws.Range("D1").Value = 0.1803224384784698
$ws.Range("E1").Value = 0.1497888714075089
